# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet, shifting the existing "Late" / heading / "Outstanding" columns one
# place to the right, and make the "Repayment schedule" sheet the active
# sheet/tab (selecting L14), while the "NewLoanInput" sheet keeps its B2
# selection but is no longer the active tab.

$wb = $excel.ActiveWorkbook

$wsNewLoanInput = $wb.Worksheets.Item("NewLoanInput")
$wsRepaymentSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column before column N (pushes N:P -> O:Q)
$wsRepaymentSchedule.Columns("N:N").Insert()

# Match the new column's stored width (11 characters, same as column M)
$wsRepaymentSchedule.Columns("N:N").ColumnWidth = 10.166666666666666

# Keep NewLoanInput's own selection as-is
[void]$wsNewLoanInput.Activate()
[void]$wsNewLoanInput.Range("B2").Select()

# Make "Repayment schedule" the active tab with the new selection
[void]$wsRepaymentSchedule.Activate()
[void]$wsRepaymentSchedule.Range("L14").Select()
